# Claim code presentations: fill in CPT Codes / HCPCS Codes / Revenue Codes
# content slides and add reference hyperlinks.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: "CPT Codes"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$s2Title = $s2.Shapes.Item(1).TextFrame.TextRange
$s2Title.InsertBefore("CPT Codes")

$s2Body = $s2.Shapes.Item(2).TextFrame.TextRange
$s2Body.InsertBefore("Current Procedural Terminology (CPT) codes identify medical procedures`rCPT codes are maintained (and copyrighted) by the American Medical Association `rThree types of CPT Codes:`rCategory 1:  commonly-used billable procedures`rCategory 2:  non-billable codes used to describe patient history, interventions, resulta and quality metrics.`rCategory 3:  Emerging technologies")

$s2Body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$s2Body2.Paragraphs(4).IndentLevel = 2
$s2Body2.Paragraphs(5).IndentLevel = 2
$s2Body2.Paragraphs(6).IndentLevel = 2

# ---------------------------------------------------------------------------
# Slide 3: "HCPCS Codes"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$s3Title = $s3.Shapes.Item(1).TextFrame.TextRange
$s3Title.InsertBefore("HCPCS Codes")

$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3Body.InsertBefore([char]0x201C + "Hick Picks" + [char]0x201D + ", the Healthcare Common Procedure coding System (HCPCS) is a coding system developed by the Center for Medicare and Medicaid Services.`rThere are three levels of HCPCS codes:`rLevel 1:  Consist of the CPT codes`rLevel 2: Records non-physician services and supplies.  These codes cover durable medical good, ambulance transportation, rehabilitation services, etc.`rLevel 3: No longer used, level 3 codes allowed state Medicaid agencies to develop custom codes.")

$s3Body2 = $s3.Shapes.Item(2).TextFrame.TextRange
$s3Body2.Paragraphs(3).IndentLevel = 2
$s3Body2.Paragraphs(4).IndentLevel = 2
$s3Body2.Paragraphs(5).IndentLevel = 2

# ---------------------------------------------------------------------------
# Slide 4: "Revenue Codes"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$s4Title = $s4.Shapes.Item(1).TextFrame.TextRange
$s4Title.InsertBefore("Revenue Codes")

$s4Shape = $s4.Shapes.Item(2)
$s4Body = $s4Shape.TextFrame.TextRange
$s4Body.InsertBefore("Revenue codes are four-digit numbers that provide additional details about a medical service`rType of room " + [char]0x2013 + " private, two-person, sterile environment, etc.`rType of medical supply " + [char]0x2013 + " sterile, non-sterile, take home supply, etc.`rType of clinic " + [char]0x2013 + " rural health clinic, family practice clinic, etc.`rThe codes are designed to distinguish between similar procedures that have different costs due to the location where the service was provided, or the supplies used to perform the procedure. ")

$s4Body2 = $s4Shape.TextFrame.TextRange
$s4Body2.Paragraphs(2).IndentLevel = 2
$s4Body2.Paragraphs(3).IndentLevel = 2
$s4Body2.Paragraphs(4).IndentLevel = 2

# Custom position/size for the content placeholder + auto-fit text
$s4Shape.Left = 120.84225
$s4Shape.Top = 158.71910
$s4Shape.Width = 749.61878
$s4Shape.Height = 316.15736
$s4Shape.TextFrame.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 5: "References" -- add two hyperlinked reference lines
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5Shape = $s5.Shapes.Item(2)
$s5Body = $s5Shape.TextFrame.TextRange

$s5Body.InsertBefore("https://c.ymcdn.com/sites/www.mpca.net/resource/resmgr/billing_and_coding/medicare%20ub-04%20revenue%20codes.pdf`rhttps://en.wikipedia.org/wiki/ICD-10`r")

$s5Body2 = $s5Shape.TextFrame.TextRange
$s5Body2.Paragraphs(1).ActionSettings.Item(1).Hyperlink.Address = "https://c.ymcdn.com/sites/www.mpca.net/resource/resmgr/billing_and_coding/medicare%20ub-04%20revenue%20codes.pdf"
$s5Body2.Paragraphs(2).ActionSettings.Item(1).Hyperlink.Address = "https://en.wikipedia.org/wiki/ICD-10"
